$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows that listed subject "حمل و نقل و زیرساخت" with
# keywords چین / ترکیه / امارات (rows 2-4), shifting the remaining rows up.
$ws.Range("A2:B4").Delete(-4162)

# Append the two new rows at the bottom of the table (now rows 28 and 29).
$ws.Range("A28").Value = "ساخت و تامین مالی"
$ws.Range("B28").Value = "شهرک"
$ws.Range("A29").Value = "مدیریت بازار"
$ws.Range("B29").Value = "معاملات"

# Update the active selection to the next empty row, matching the saved file.
$ws.Range("A30").Select()
